$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.595.24"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.925.24"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.39"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08232"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.84"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.924.57"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.110"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.276"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.91"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06887"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.63"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.601.28"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.685"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.03"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.190"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.178.18"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.14"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.436"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.093"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.76"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.014"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09637"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.625"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.578"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.378"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06390"
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02292"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.189"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5954"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.74"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.872"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1853"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.427"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.247"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.39"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5566"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.979"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.62"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.443"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.14"
$ws.Range("E51").Value = "  -0.43%  "
